# feat: add 2022-Q3 data
#
# 1. Update the "总计" (total) summary sheet: insert a new leading data
#    row for "2022-Q3" (count=1, value=0.3) and keep the rest of the
#    quarters, with the running index column (A) renumbered 0..6.
# 2. Insert a brand-new worksheet named "2022-Q3", placed right after
#    "总计" (so it becomes the first quarterly detail sheet), and fill
#    it in with the same layout/styling as the other quarterly sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$quarters = @(
    @{ Label = "2022-Q3"; Count = 1;  Value = 0.3 },
    @{ Label = "2022-Q2"; Count = 7;  Value = 0.25 },
    @{ Label = "2022-Q1"; Count = 4;  Value = 0.39 },
    @{ Label = "2021-Q4"; Count = 24; Value = 1.94 },
    @{ Label = "2021-Q3"; Count = 5;  Value = 0.29 },
    @{ Label = "2021-Q2"; Count = 2;  Value = 0.06 },
    @{ Label = "2021-Q1"; Count = 2;  Value = 0.05 }
)

for ($i = 0; $i -lt $quarters.Count; $i++) {
    $r = $i + 2
    $q = $quarters[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $q.Label
    $total.Cells.Item($r, 3).Value = $q.Count
    $total.Cells.Item($r, 4).Value = $q.Value
}

# Row 8 is brand new (the sheet previously ended at row 7), so column A's
# bold/centered/bordered "index" style has to be applied explicitly -
# rows 2-7 already carry it from the pre-existing cells.
$lastIndexCell = $total.Cells.Item($quarters.Count + 1, 1)
$lastIndexCell.Font.Bold = $true
$lastIndexCell.HorizontalAlignment = -4108
$lastIndexCell.VerticalAlignment = -4160
$lastIndexCell.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2) New "2022-Q3" detail sheet, inserted right after "总计"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$new = $wb.Worksheets.Add($null, $afterSheet)
$new.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $new.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$new.Cells.Item(2, 1).Value = 0

# Header row + running-index column share the same bold/centered/bordered
# style used throughout the workbook's other quarterly sheets.
$headerStyle = $new.Range("B1:H1")
$headerStyle.Font.Bold = $true
$headerStyle.HorizontalAlignment = -4108
$headerStyle.VerticalAlignment = -4160
$headerStyle.Borders.LineStyle = 1

$indexStyle = $new.Range("A2")
$indexStyle.Font.Bold = $true
$indexStyle.HorizontalAlignment = -4108
$indexStyle.VerticalAlignment = -4160
$indexStyle.Borders.LineStyle = 1

$new.Cells.Item(2, 2).NumberFormat = "@"
$new.Cells.Item(2, 2).Value = "007835"

$new.Cells.Item(2, 3).NumberFormat = "@"
$new.Cells.Item(2, 3).Value = "国泰鑫睿混合"

$new.Cells.Item(2, 4).NumberFormat = "@"
$new.Cells.Item(2, 4).Value = "8.30"

$new.Cells.Item(2, 5).NumberFormat = "@"
$new.Cells.Item(2, 5).Value = "79.49"

$new.Cells.Item(2, 6).NumberFormat = "@"
$new.Cells.Item(2, 6).Value = "3.62"

$new.Cells.Item(2, 7).NumberFormat = "@"
$new.Cells.Item(2, 7).Value = "0.3005"

$new.Cells.Item(2, 8).Value = 5
